# Jogos_do_Dia_Betfair_Back_Lay_2026-02-27.xlsx update
# - refresh the snapshot timestamp (column BH) on every row
# - update odds for several existing fixtures (rows 4, 6, 7)
# - turn row 8 into the Colombian Primera B / Bogota vs Atletico FC Cali
#   fixture (odds columns F:BG stay as they were)
# - insert two brand-new fixtures as rows 9 and 10 (Colombian Primera B /
#   Leones FC vs Independiente Yumbo, and Algerian Ligue 1 / JS Saoura vs
#   ES Ben Aknoun — i.e. what used to live in row 8)
# - push the Uruguayan Primera Division fixture from row 9 down to row 11,
#   refreshing its odds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTs = "2026-02-25 03:16:49"

# ---------------------------------------------------------------------------
# 1) Insert two blank rows at 9 and 10, pushing the old row 9 (Uruguayan
#    Primera Division / Cerro vs Boston River) down to row 11.
# ---------------------------------------------------------------------------
$ws.Rows("9:10").Insert()

# ---------------------------------------------------------------------------
# 2) Refresh the snapshot timestamp (column BH) on every data row.
# ---------------------------------------------------------------------------
foreach ($r in 2..11) {
    $ws.Range("BH$r").Value = $newTs
}

# ---------------------------------------------------------------------------
# 3) Row 4 - Al-Shabab (KSA) vs Al-Hilal - updated odds.
# ---------------------------------------------------------------------------
$row4 = @{
    "F4" = 1.04
    "G4" = 14
    "H4" = 1.41
    "I4" = 1.48
    "P4" = 2.28
    "Q4" = 1.56
}
foreach ($key in $row4.Keys) {
    $ws.Range($key).Value = $row4[$key]
}

# ---------------------------------------------------------------------------
# 4) Row 6 - Al-Ittihad vs Al-Khaleej Saihat - updated odds.
# ---------------------------------------------------------------------------
$row6 = @{
    "F6" = 1.51
    "G6" = 1.63
    "H6" = 4.8
    "I6" = 7.2
    "J6" = 4.6
    "K6" = 6.6
    "P6" = 2.8
    "Q6" = 1.46
}
foreach ($key in $row6.Keys) {
    $ws.Range($key).Value = $row6[$key]
}

# ---------------------------------------------------------------------------
# 5) Row 7 - Albacete vs Almeria - updated odds.
# ---------------------------------------------------------------------------
$row7 = @{
    "F7" = 2.46
    "G7" = 2.86
    "H7" = 2.78
    "I7" = 3.3
    "J7" = 3.45
    "K7" = 4.1
    "P7" = 2.06
    "Q7" = 1.64
}
foreach ($key in $row7.Keys) {
    $ws.Range($key).Value = $row7[$key]
}

# ---------------------------------------------------------------------------
# 6) Row 8 - was Algerian Ligue 1 / JS Saoura vs ES Ben Aknoun, now becomes
#    Colombian Primera B / Bogota vs Atletico FC Cali. Odds columns (F:BG)
#    are left untouched. Column B ("2026-02-27") is untouched too.
# ---------------------------------------------------------------------------
$row8 = @{
    "A8" = "Colombian Primera B"
    "C8" = "17:00:00"
    "D8" = "Bogota"
    "E8" = "Atletico FC Cali"
}
foreach ($key in $row8.Keys) {
    $ws.Range($key).Value = $row8[$key]
}

# ---------------------------------------------------------------------------
# 7) New row 9 - Colombian Primera B / Leones FC vs Independiente Yumbo.
#    New row 10 - Algerian Ligue 1 / JS Saoura vs ES Ben Aknoun.
#
#    Column B just needs the literal text "2026-02-27" (same date as every
#    other fixture). Assigning that string straight to .Value makes the COM
#    layer auto-convert it into a date serial, so instead we copy it out of
#    an existing cell that already holds it as plain text (B2) — Range.Copy
#    carries the text value across without reinterpreting it.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("B9"))
$ws.Range("B2").Copy($ws.Range("B10"))

$row9 = @{
    "A9" = "Colombian Primera B"
    "C9" = "17:00:00"
    "D9" = "Leones FC"
    "E9" = "Independiente Yumbo"
}
foreach ($key in $row9.Keys) {
    $ws.Range($key).Value = $row9[$key]
}

$row10 = @{
    "A10" = "Algerian Ligue 1"
    "C10" = "18:00:00"
    "D10" = "JS Saoura"
    "E10" = "ES Ben Aknoun"
}
foreach ($key in $row10.Keys) {
    $ws.Range($key).Value = $row10[$key]
}

# Rows 9 and 10 share the same numeric (F:BG) pattern that row 8 had before
# this edit: F=1.04 G=1000 H=1.04 I=1000 J=1.01 K=950, L-O=0, P=1.24 Q=1.01,
# R-BG=0.
$numericCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z",
    "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT",
    "AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG")

$sharedPattern = @{
    "F" = 1.04
    "G" = 1000
    "H" = 1.04
    "I" = 1000
    "J" = 1.01
    "K" = 950
    "P" = 1.24
    "Q" = 1.01
}

foreach ($col in $numericCols) {
    if ($sharedPattern.ContainsKey($col)) {
        $val = $sharedPattern[$col]
    } else {
        $val = 0
    }
    $ws.Range("${col}9").Value = $val
    $ws.Range("${col}10").Value = $val
}

# ---------------------------------------------------------------------------
# 8) Row 11 (previously row 9) - Uruguayan Primera Division / Cerro vs
#    Boston River - updated odds.
# ---------------------------------------------------------------------------
$row11 = @{
    "F11" = 2.98
    "G11" = 4
    "H11" = 2.46
    "I11" = 2.82
    "J11" = 2.8
    "K11" = 3.55
    "P11" = 1.48
    "Q11" = 2.62
}
foreach ($key in $row11.Keys) {
    $ws.Range($key).Value = $row11[$key]
}
